{"js": "// Add a new paragraph \"Hey dude, it's been a long time since we met.\"\n// right after the last existing paragraph in the document body.\nconst body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\n\nlastParagraph.insertParagraph(\n  \"Hey dude, it\\u2019s been a long time since we met.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Add a new paragraph \"Hey dude, it's been a long time since we met.\"\n# right after the last existing paragraph in the document body.\n$d = $word.ActiveDocument\n\n$lastRange = $d.Paragraphs.Last.Range\n$lastRange.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.Text = \"Hey dude, it\u2019s been a long time since we met.\"\n"}
